# chore: persist data folder changes (clientes, historial, docs)
#
# Adds the new cliente record C1004 / SAUL TORRES as the next row in the
# "Clientes" sheet (data/clientes.xlsx), mirroring the columns already
# populated for the other rows (id, nombre, sucursal, asesor,
# fecha_ingreso, fecha_dispersion, estatus).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row right after the existing data (row 5 -> row 6).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "C1004"
$ws.Cells.Item($newRow, 2).Value = "SAUL TORRES"
$ws.Cells.Item($newRow, 3).Value = "TOXQUI"
$ws.Cells.Item($newRow, 4).Value = "Martha Ortiz"

# fecha_ingreso / fecha_dispersion are stored as plain text (e.g. "2025-10-08")
# in this sheet, not real dates, so force Text formatting before writing the
# value and then drop the formatting again so the cell ends up styled just
# like its neighbours (no explicit number format left behind).
$dateRange = $ws.Range($ws.Cells.Item($newRow, 5), $ws.Cells.Item($newRow, 6))
$dateRange.NumberFormat = "@"
$ws.Cells.Item($newRow, 5).Value = "2025-10-08"
$ws.Cells.Item($newRow, 6).Value = "2025-10-08"
$dateRange.ClearFormats()

$ws.Cells.Item($newRow, 7).Value = "DISPERSADO"
